# This script reproduces the NATMI LR-pair data refresh for Cort-Sstr2:
#  - sending/target cluster labels are rewritten (new label set / order)
#  - numeric score columns are updated with newly computed TPM-based values
#  - the table grows from 12 data rows (rows 2-13) to 15 data rows (rows 2-16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old data block (header row 1 is left untouched) so the shared-string
# table and row count get rebuilt cleanly for the new data.
$ws.Range("A2:T13").Clear()

# New data rows (Sending cluster, Ligand, Receptor, Target cluster, then the 16 score columns)
$data = @(
  @("FAPs", "Cort", "Sstr2", "ECs", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 1, 0.3333333333333333, 0.213999, 0.641997, 0.5677164163662266, 0.5677164163662265, 0.103044869813, 0.9274038283170002, 0.442250776077402, 0.4422507760774019),
  @("FAPs", "Cort", "Sstr2", "FAPs", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 1, 0.3333333333333333, 0.02487466666666667, 0.074624, 0.06598982527163412, 0.06598982527163411, 0.01197765778488889, 0.107798920064, 0.05140603758896076, 0.05140603758896076),
  @("FAPs", "Cort", "Sstr2", "Inflammatory-Mac", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 1, 0.3333333333333333, 0.04845433333333333, 0.145363, 0.1285441543063967, 0.1285441543063967, 0.02333174673811111, 0.209985720643, 0.1001358254990901, 0.1001358254990901),
  @("FAPs", "Cort", "Sstr2", "MuSCs", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 1, 0.3333333333333333, 0.01566266666666667, 0.046988, 0.04155137636502391, 0.0415513763650239, 0.007541892474222224, 0.06787703226800001, 0.03236849933305758, 0.03236849933305757),
  @("FAPs", "Cort", "Sstr2", "Resolving-Mac", 2, 0.6666666666666666, 0.4815203333333334, 1.444561, 0.7789994499509271, 0.7789994499509271, 1, 0.3333333333333333, 0.07395633333333333, 0.221869, 0.1961982276907187, 0.1961982276907187, 0.03561147827877778, 0.3205033045090001, 0.1528383114524166, 0.1528383114524166),
  @("MuSCs", "Cort", "Sstr2", "ECs", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 1, 0.3333333333333333, 0.213999, 0.641997, 0.5677164163662266, 0.5677164163662265, 0.011747475105, 0.105727275945, 0.05041813330022542, 0.05041813330022541),
  @("MuSCs", "Cort", "Sstr2", "FAPs", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 1, 0.3333333333333333, 0.02487466666666667, 0.074624, 0.06598982527163412, 0.06598982527163411, 0.001365494826666667, 0.01228945344, 0.005860467851712736, 0.005860467851712734),
  @("MuSCs", "Cort", "Sstr2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 1, 0.3333333333333333, 0.04845433333333333, 0.145363, 0.1285441543063967, 0.1285441543063967, 0.002659900628333333, 0.023939105655, 0.01141583389162358, 0.01141583389162358),
  @("MuSCs", "Cort", "Sstr2", "MuSCs", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 1, 0.3333333333333333, 0.01566266666666667, 0.046988, 0.04155137636502391, 0.0415513763650239, 0.0008598020866666668, 0.00773821878, 0.003690121990462559, 0.003690121990462559),
  @("MuSCs", "Cort", "Sstr2", "Resolving-Mac", 1, 0.3333333333333333, 0.054895, 0.164685, 0.08880865841952565, 0.08880865841952565, 1, 0.3333333333333333, 0.07395633333333333, 0.221869, 0.1961982276907187, 0.1961982276907187, 0.004059832918333333, 0.036538496265, 0.01742410138550135, 0.01742410138550135),
  @("Resolving-Mac", "Cort", "Sstr2", "ECs", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 1, 0.3333333333333333, 0.213999, 0.641997, 0.5677164163662266, 0.5677164163662265, 0.017486143622, 0.157375292598, 0.0750475069885992, 0.07504750698859919),
  @("Resolving-Mac", "Cort", "Sstr2", "FAPs", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 1, 0.3333333333333333, 0.02487466666666667, 0.074624, 0.06598982527163412, 0.06598982527163411, 0.002032542179555555, 0.018292879616, 0.008723319830960622, 0.00872331983096062),
  @("Resolving-Mac", "Cort", "Sstr2", "Inflammatory-Mac", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 1, 0.3333333333333333, 0.04845433333333333, 0.145363, 0.1285441543063967, 0.1285441543063967, 0.003959268182444444, 0.035633413642, 0.01699249491568301, 0.01699249491568301),
  @("Resolving-Mac", "Cort", "Sstr2", "MuSCs", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 1, 0.3333333333333333, 0.01566266666666667, 0.046988, 0.04155137636502391, 0.0415513763650239, 0.001279817376888889, 0.011518356392, 0.005492755041503775, 0.005492755041503774),
  @("Resolving-Mac", "Cort", "Sstr2", "Resolving-Mac", 1, 0.3333333333333333, 0.08171133333333333, 0.245134, 0.1321918916295473, 0.1321918916295473, 1, 0.3333333333333333, 0.07395633333333333, 0.221869, 0.1961982276907187, 0.1961982276907187, 0.006043070605111111, 0.054387635446, 0.02593581485280073, 0.02593581485280073)
)

$r = 2
foreach ($row in $data) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value2 = $val
    $c++
  }
  $r++
}

